{"js": "// Fix the two \"lost closing tag\" typos in the mustache-style template\n// text that lives inside the Word body:\n//   1. \"{{reportTime}\"  ->  \"{{reportTime}}\"      (missing closing brace)\n//   2. \"`<body>...`\"    ->  \"`<meta charset=\"UTF-8\"><body>...`\"\n//      (missing <meta charset=\"UTF-8\"> tag right before the <body> tag)\n//\n// Both fixes are plain text insertions right next to an unambiguous,\n// unique anchor string, so `Body.search` + `Range.insertText` is enough -\n// no Find/Replace wildcards or content-control gymnastics needed.\n\nconst body = context.document.body;\n\n// 1) \"{{reportTime}\" -> \"{{reportTime}}\"\n// Insert a new \"}\" immediately after the existing \"reportTime\" text\n// (i.e. just before the lone \"}\" that is already there), producing the\n// doubled closing brace the mustache syntax needs.\nconst reportTimeResults = body.search(\"reportTime\", { matchCase: true });\nreportTimeResults.load(\"items\");\nawait context.sync();\n\nif (reportTimeResults.items.length > 0) {\n  reportTimeResults.items[0].insertText(\"}\", Word.InsertLocation.after);\n}\n\n// 2) Insert `<meta charset=\"UTF-8\">` right before the existing `<body>`\n// literal text (both live inside the same template-literal run, right\n// after the opening back-tick).\nconst bodyTagResults = body.search(\"<body>\", { matchCase: true });\nbodyTagResults.load(\"items\");\nawait context.sync();\n\nif (bodyTagResults.items.length > 0) {\n  bodyTagResults.items[0].insertText('<meta charset=\"UTF-8\">', Word.InsertLocation.before);\n}\n\nawait context.sync();\n", "ps1": "# Fix the two \"lost closing tag\" typos in the mustache-style template\n# text that lives inside the Word body:\n#   1. \"{{reportTime}\"  ->  \"{{reportTime}}\"      (missing closing brace)\n#   2. \"`<body>...`\"    ->  \"`<meta charset=\"UTF-8\"><body>...`\"\n#      (missing <meta charset=\"UTF-8\"> tag right before the <body> tag)\n#\n# Both fixes are plain text insertions right next to an unambiguous,\n# unique anchor string, so Range.Find.Execute + InsertAfter/InsertBefore\n# is all that's needed.\n\n$d = $word.ActiveDocument\n\n# 1) \"{{reportTime}\" -> \"{{reportTime}}\"\n# Find \"reportTime\", collapse the found range to its end point (right\n# before the lone \"}\" that is already there) and insert a new \"}\",\n# producing the doubled closing brace the mustache syntax needs.\n$rngReportTime = $d.Content\n$found1 = $rngReportTime.Find.Execute(\"reportTime\")\nif ($found1) {\n    $rngReportTime.Collapse(0)   # wdCollapseEnd\n    $rngReportTime.InsertAfter(\"}\")\n}\n\n# 2) Insert `<meta charset=\"UTF-8\">` right before the existing `<body>`\n# literal text (both live inside the same template-literal run, right\n# after the opening back-tick).\n$rngBody = $d.Content\n$found2 = $rngBody.Find.Execute(\"<body>\")\nif ($found2) {\n    $rngBody.Collapse(1)   # wdCollapseStart\n    $rngBody.InsertBefore('<meta charset=\"UTF-8\">')\n}\n"}
